$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 325, shifting existing rows 325:393 down to 326:394
$ws.Rows(325).Insert()

# Populate the newly inserted row 325 with the new record
$ws.Cells.Item(325, 1).Value = 10
$ws.Cells.Item(325, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(325, 3).Value = "La Araucanía"
$ws.Cells.Item(325, 4).Value = 44995
$ws.Cells.Item(325, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(325, 5).Value = 9
$ws.Cells.Item(325, 6).Value = "Fruta"
$ws.Cells.Item(325, 7).Value = 100102
$ws.Cells.Item(325, 8).Value = "Cítricos"
$ws.Cells.Item(325, 9).Value = 100102006
$ws.Cells.Item(325, 10).Value = "Pomelo"
$ws.Cells.Item(325, 11).Value = "Start Ruby"
$ws.Cells.Item(325, 12).Value = "Primera"
$ws.Cells.Item(325, 13).Value = 100
$ws.Cells.Item(325, 14).Value = 14000
$ws.Cells.Item(325, 15).Value = 14000
$ws.Cells.Item(325, 16).Value = 14000
$ws.Cells.Item(325, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(325, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(325, 19).Value = 933
$ws.Cells.Item(325, 20).Value = 15
